$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume (E) columns to remain text so that
# numeric-looking strings (e.g. "1.00", "563.41") are not auto-converted
# to numbers by Excel when the .Value is assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.451.55'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '2.991.89'
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '563.41'
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").Value = '139.16'
$ws.Range("E6").Value = '  +5.50%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.61%  '
$ws.Range("D9").Value = '2.983.54'
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("E10").Value = '  +4.31%  '
$ws.Range("D11").Value = '5.31'
$ws.Range("E11").Value = '  +11.66%  '
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("D13").Value = '0.0000230'
$ws.Range("E13").Value = '  +4.44%  '
$ws.Range("D14").Value = '33.81'
$ws.Range("E14").Value = '  +3.14%  '
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '3.489.72'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").Value = '7.17'
$ws.Range("E17").Value = '  +4.67%  '
$ws.Range("D18").Value = '2.991.19'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").Value = '59.459.60'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("D20").Value = '435.11'
$ws.Range("E20").Value = '  +4.36%  '
$ws.Range("D21").Value = '13.63'
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").Value = '0.718'
$ws.Range("E22").Value = '  +3.57%  '
$ws.Range("D23").Value = '13.40'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = '7.03'
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").Value = '80.02'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("E27").Value = '  +10.34%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("E29").Value = '  +3.07%  '
$ws.Range("D30").Value = '7.78'
$ws.Range("E30").Value = '  +5.07%  '
$ws.Range("D31").Value = '6.25'
$ws.Range("E31").Value = '  +4.69%  '
$ws.Range("E32").Value = '  +9.18%  '
$ws.Range("D33").Value = '25.75'
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("D34").Value = '0.0₃0778'
$ws.Range("E34").Value = '  +12.15%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +6.43%  '
$ws.Range("D36").Value = '5.89'
$ws.Range("E36").Value = '  +3.74%  '
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").Value = '48.83'
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").Value = '8.61'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("D40").Value = '2.78'
$ws.Range("E40").Value = '  +6.69%  '
$ws.Range("D41").Value = '400.93'
$ws.Range("E41").Value = '  +7.12%  '
$ws.Range("D42").Value = '0.0354'
$ws.Range("E42").Value = '  +2.68%  '
$ws.Range("D43").Value = '2.759.76'
$ws.Range("E43").Value = '  +2.65%  '
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("E45").Value = '  +6.30%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").Value = '34.65'
$ws.Range("E47").Value = '  +18.51%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '122.84'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("E49").Value = '  +1.87%  '
$ws.Range("E50").Value = '  +2.96%  '
$ws.Range("D51").Value = '23.52'
$ws.Range("E51").Value = '  +2.72%  '
